# Add a new "DataDriven" worksheet after the last sheet (PNR), populate it
# with 5 rows of the (username, password) credential pair used elsewhere in
# this workbook ("bakshu405" / "Bakshu405"), auto-fit its columns, and make
# it the active sheet/tab - matching the IRCTC XLData.xlsx commit that
# introduced a DataDriven test-data sheet.

$wb = $excel.ActiveWorkbook

# New sheet goes after the current last worksheet (PNR), so it lands at the
# end of the tab strip.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "DataDriven"

# Fill A1:B5 with the username/password pair, 5 rows deep.
for ($row = 1; $row -le 5; $row++) {
    $newSheet.Cells.Item($row, 1).Value = "bakshu405"
    $newSheet.Cells.Item($row, 2).Value = "Bakshu405"
}

# Auto-size the two data columns to fit their contents.
$newSheet.Range("A1:B5").EntireColumn.AutoFit() | Out-Null

# Leave the cursor on A5 - this also leaves "DataDriven" as the
# active/selected sheet, mirroring the workbook's final saved view state.
[void]$newSheet.Range("A5").Select()
